$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 270
$ws.Range("B3").Value = 262
$ws.Range("B4").Value = 283
$ws.Range("B5").Value = 291
$ws.Range("B6").Value = 295
$ws.Range("B7").Value = 285
$ws.Range("B8").Value = 275
$ws.Range("B9").Value = 299
$ws.Range("B10").Value = 283
$ws.Range("B11").Value = 298
$ws.Range("B12").Value = 292
$ws.Range("B13").Value = 289

$ws.Range("B13").Select()
